$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All source cells in this sheet are stored as text (inline strings),
# so force text NumberFormat before assigning to avoid Excel auto-converting
# numeric-looking strings (e.g. "546.78", "2.976.61") into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.179.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.976.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.45"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.92%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.969.32"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.89"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -7.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000218"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.51"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.463.69"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.240.71"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.109"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.976.21"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.53"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.76"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.660"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.89"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.87"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.54"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.85%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.87"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.26"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.28"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.43"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "54.30"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.80"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "442.80"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -10.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.124.90"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0783"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0376"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.01%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -11.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.24"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.239"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.40%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.92"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.37%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "BitgetToken"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +9.64%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.58"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₃0477"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -10.42%  "
